# Trade #39 closed at 2026-02-17 21:03:29 - unknown UNKNOWN +0.000%
#
# This script updates the "live_trading_results" workbook to reflect:
#  1) Updated aggregate metrics on the Summary sheet
#  2) Updated aggregate metrics for the MarketMaking strategy on Strategy Status
#  3) Trade #67 (row 68 on "All Trades", row 35 on "MarketMaking") closing out
#     with an early_exit, plus a brand-new open Trade #100 appended to both
#     the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a literal text value into a cell without Excel's automatic
# date/number recognition kicking in (e.g. "2026-02-17" turning into a date
# serial). We flip the cell to Text format, assign the literal string, then
# clear the formatting again so no stray style survives on the cell.
# ---------------------------------------------------------------------------
function Set-TextValue($ws, $addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

# ===========================================================================
# 1) Summary sheet
# ===========================================================================
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.58
$summary.Range("B4").Value = 0.38
$summary.Range("B6").Value = 67
$summary.Range("B7").Value = 32
$summary.Range("B9").Value = 47.76

# ===========================================================================
# 2) Strategy Status sheet (MarketMaking row)
# ===========================================================================
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.58
$status.Range("D5").Value = 34
$status.Range("E5").Value = 0.27
$status.Range("F5").Value = 0.58
$status.Range("G5").Value = 52.94

# ===========================================================================
# 3) All Trades sheet
# ===========================================================================
$allTrades = $wb.Worksheets.Item("All Trades")

# --- Trade #67 (row 68) now closes out early ---
$allTrades.Range("G68").Value = 0.063846
$allTrades.Range("H68").Value = "CLOSED"
$allTrades.Range("I68").Value = 112.8216
$allTrades.Range("J68").Value = 0.03
$allTrades.Range("K68").Value = 100.58
Set-TextValue $allTrades "L68" "early_exit"
$allTrades.Range("M68").Value = 0.13

# --- New Trade #100 (row 101) appended ---
$allTrades.Range("A101").Value = 100
Set-TextValue $allTrades "B101" "2026-02-17"
Set-TextValue $allTrades "C101" "21:03:22"
Set-TextValue $allTrades "D101" "MarketMaking"
Set-TextValue $allTrades "E101" "UP"
$allTrades.Range("F101").Value = 0.03
$allTrades.Range("H101").Value = "OPEN"
$allTrades.Range("I101").Value = 0
$allTrades.Range("J101").Value = 0
$allTrades.Range("K101").Value = 100.5510412885904
$allTrades.Range("M101").Value = 0
$allTrades.Range("N101").Value = 0
$allTrades.Range("O101").Value = 0
$allTrades.Range("P101").Value = 0.6
Set-TextValue $allTrades "Q101" "Normal spread capture: 19600 bps"

# ===========================================================================
# 4) MarketMaking sheet
# ===========================================================================
$mm = $wb.Worksheets.Item("MarketMaking")

# --- Trade #67 (row 35) now closes out early ---
$mm.Range("G35").Value = 0.063846
$mm.Range("H35").Value = "CLOSED"
$mm.Range("I35").Value = 112.8216
$mm.Range("J35").Value = 0.03
$mm.Range("K35").Value = 100.58
Set-TextValue $mm "P35" "early_exit"
$mm.Range("Q35").Value = 0.13

# --- New Trade #100 (row 68) appended ---
$mm.Range("A68").Value = 100
Set-TextValue $mm "B68" "2026-02-17"
Set-TextValue $mm "C68" "21:03:22"
Set-TextValue $mm "D68" "MarketMaking"
Set-TextValue $mm "E68" "UP"
$mm.Range("F68").Value = 0.03
$mm.Range("H68").Value = "OPEN"
$mm.Range("I68").Value = 0
$mm.Range("J68").Value = 0
$mm.Range("K68").Value = 100.5510412885904
$mm.Range("L68").Value = 0
$mm.Range("M68").Value = 0
$mm.Range("N68").Value = 0.6
Set-TextValue $mm "O68" "Normal spread capture: 19600 bps"
$mm.Range("Q68").Value = 0

Write-Output "edit complete"
